# Update and add file Acme System
# Adds the new "Calculate Client Security Hash" asset configuration to the
# Settings sheet, bumps the MaxRetryNumber constant, and leaves the
# workbook positioned on the Constants tab (matching the author's final
# selection state).

$wb = $excel.ActiveWorkbook

# --- Settings sheet -------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# New asset rows (System1_URL, SHA1Online_URL, acme-test_Credential).
# Filled in this order so new shared-string entries land in the same
# sequence as the authored workbook (names, then values, then credential).
$settings.Range("A6").Value = "System1_URL"
$settings.Range("A7").Value = "SHA1Online_URL"
$settings.Range("B6").Value = "https://acme-test.uipath.com/login"
$settings.Range("B7").Value = "http://www.sha1-online.com/"
$settings.Range("A8").Value = "acme-test_Credential"
$settings.Range("B8").Value = "acme-test_Credential"

# Business process name used for logging now reflects the new process.
$settings.Range("B5").Value = "Calculate Client Security Hash"

$settings.Range("B6").Select()

# --- Constants sheet --------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

# MaxRetryNumber raised from 0 to 2.
$constants.Range("B2").Value = 2

$constants.Activate()
$constants.Range("B3").Select()
